# The "Reading and Writing Data" mini-section (originally 4 consecutive
# slides: intro, "Common File Formats", "Reading File Formats",
# "Writing File Formats") is moved earlier in the deck, ahead of the
# "Data Manipulation" / "Subsetting" / "Filtering" slides, so that file
# I/O is covered before the data-manipulation topics that build on it.
#
# Before: ... 19) TYPE CONVERSION, 20) DATA MANIPULATION, 21) SUBSETTING,
#          22) FILTERING, 23) READING AND WRITING DATA, 24) COMMON FILE
#          FORMATS, 25) READING FILE FORMATS, 26) WRITING FILE FORMATS,
#          27) VISUALIZATION ...
# After:  ... 19) TYPE CONVERSION, 20) READING AND WRITING DATA,
#          21) COMMON FILE FORMATS, 22) READING FILE FORMATS,
#          23) WRITING FILE FORMATS, 24) DATA MANIPULATION,
#          25) SUBSETTING, 26) FILTERING, 27) VISUALIZATION ...

$p = $ppt.ActivePresentation

# Move the "READING AND WRITING DATA" intro slide (currently #23) so it
# lands right before the old #20 ("DATA MANIPULATION"), i.e. becomes #20.
$p.Slides.Item(23).MoveTo(20)

# Move "COMMON FILE FORMATS" (now at #24) to right after it, as #21.
$p.Slides.Item(24).MoveTo(21)

# Move "READING FILE FORMATS" (now at #25) to right after that, as #22.
$p.Slides.Item(25).MoveTo(22)

# Move "WRITING FILE FORMATS" (now at #26) to right after that, as #23.
$p.Slides.Item(26).MoveTo(23)
